# Auto-applies the cell-level changes described by the source diff.
# Numeric ("want-to-go" counts) and text (title/venue/time/link/cover) updates
# across the "展览", "演出" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# ---- Sheet: 展览 ----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 4594
$ws.Range("F3").Value = 2698
$ws.Range("F5").Value = 2693
$ws.Range("F9").Value = 1702
$ws.Range("F10").Value = 723
$ws.Range("F12").Value = 176
$ws.Range("F13").Value = 377
$ws.Range("F15").Value = 287
$ws.Range("F16").Value = 84
$ws.Range("F17").Value = 59
$ws.Range("F18").Value = 507
$ws.Range("F21").Value = 630
$ws.Range("F22").Value = 729
$ws.Range("F24").Value = 21
$ws.Range("F25").Value = 483
$ws.Range("F26").Value = 1648
$ws.Range("F27").Value = 1375
$ws.Range("F28").Value = 283
$ws.Range("F29").Value = 35
$ws.Range("F30").Value = 1362
$ws.Range("F31").Value = 2222
$ws.Range("F32").Value = 356
$ws.Range("F34").Value = 585
$ws.Range("F38").Value = 745
$ws.Range("F39").Value = 1425
$ws.Range("F40").Value = 177
$ws.Range("F42").Value = 468
$ws.Range("F44").Value = 103

# ---- Sheet: 演出 ----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 4

# ---- Sheet: 全部类型 ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 4594
$ws.Range("F3").Value = 2698
$ws.Range("F4").Value = 2693
$ws.Range("F5").Value = 1702
$ws.Range("F7").Value = 4
$ws.Range("F8").Value = 723
$ws.Range("F10").Value = 176
$ws.Range("F11").Value = 377
$ws.Range("F13").Value = 287
$ws.Range("F14").Value = 84
$ws.Range("F15").Value = 59
$ws.Range("F16").Value = 507
$ws.Range("F18").Value = 630
$ws.Range("F19").Value = 729
$ws.Range("F24").Value = 21
$ws.Range("F25").Value = 483
$ws.Range("F26").Value = 1648
$ws.Range("F27").Value = 1375
$ws.Range("F28").Value = 283
$ws.Range("F29").Value = 35
$ws.Range("C31").Value = "杭州·“真的爱你”致敬Beyond·黄家驹31周年演唱会-630乐团再现91殿堂级演出"
$ws.Range("D31").Value = "中山南路77号尚城·利星1157 3F MAOLivehouse杭州"
$ws.Range("E31").Value = "2024.05.02 20:00-05.02 22:00"
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 98
$ws.Range("H31").Value = "https://show.bilibili.com/platform/detail.html?id=83545"
$ws.Range("I31").Value = "//i2.hdslb.com/bfs/openplatform/202403/1VblWbtA1711696632442.jpeg"
$ws.Range("C32").Value = "杭州·第四届华盟动漫次元嘉年华"
$ws.Range("D32").Value = "创意路1号 中国智谷富春园区"
$ws.Range("E32").Value = "2024.05.02 10:00-05.03 17:00"
$ws.Range("F32").Value = 2222
$ws.Range("G32").Value = 58
$ws.Range("H32").Value = "https://show.bilibili.com/platform/detail.html?id=82465"
$ws.Range("I32").Value = "//i0.hdslb.com/bfs/openplatform/202403/4XHyqi3D1709780326858.jpeg"
$ws.Range("C33").Value = "杭州·造梦探险家Porject6野蛮冲撞——第五人格ONLY"
$ws.Range("D33").Value = "欢西路1号 天都城酒店"
$ws.Range("E33").Value = "2024.05.02 10:00-05.02 22:00"
$ws.Range("F33").Value = 356
$ws.Range("G33").Value = 28
$ws.Range("H33").Value = "https://show.bilibili.com/platform/detail.html?id=82851"
$ws.Range("I33").Value = "//i1.hdslb.com/bfs/openplatform/202403/a7IYN66u1711441126355.png"
$cell = $ws.Range("B34")
$cell.NumberFormat = "@"
$cell.Value = "2024-05-12"
$cell.Style = "Normal"
$ws.Range("C34").Value = "杭州·《卡农》永恒经典名曲音乐会"
$ws.Range("D34").Value = "武林路77号 浙江省文化馆小剧场（原群艺馆小剧场）"
$ws.Range("E34").Value = "2024.05.12 14:00-05.12 15:30"
$ws.Range("F34").Value = 2
$ws.Range("G34").Value = 100
$ws.Range("H34").Value = "https://show.bilibili.com/platform/detail.html?id=83176"
$ws.Range("I34").Value = "//i0.hdslb.com/bfs/openplatform/202403/gLrSkh0O1711013683966.jpeg"
$ws.Range("C35").Value = "杭州·奇迹の闪耀 「UP!」巡回动漫演唱会"
$ws.Range("D35").Value = "东坡路10号 杭州东坡大剧院"
$ws.Range("E35").Value = "2024.05.12 19:30-05.12 21:30"
$ws.Range("F35").Value = 12
$ws.Range("G35").Value = 126
$ws.Range("H35").Value = "https://show.bilibili.com/platform/detail.html?id=82452"
$ws.Range("I35").Value = "//i1.hdslb.com/bfs/openplatform/202403/HvxHPz981709707512970.jpeg"
$cell = $ws.Range("B36")
$cell.NumberFormat = "@"
$cell.Value = "2024-05-18"
$cell.Style = "Normal"
$ws.Range("C36").Value = "杭州·Jo迪"
$ws.Range("D36").Value = "萧杭路28号 格拉斯club"
$ws.Range("E36").Value = "2024.05.18 13:00-05.18 19:00"
$ws.Range("F36").Value = 19
$ws.Range("G36").Value = 198
$ws.Range("H36").Value = "https://show.bilibili.com/platform/detail.html?id=83008"
$ws.Range("I36").Value = "//i1.hdslb.com/bfs/openplatform/202403/AEtl5BHN1711015003341.jpeg"
$ws.Range("C37").Value = "杭州·《沐云华·次元狂想》经典动漫二次元ACG音乐会"
$ws.Range("D37").Value = "建国南路280号 杭州红星剧院"
$ws.Range("E37").Value = "2024.05.18 19:30-05.18 22:00"
$ws.Range("F37").Value = 15
$ws.Range("G37").Value = 90
$ws.Range("H37").Value = "https://show.bilibili.com/platform/detail.html?id=83113"
$ws.Range("I37").Value = "//i1.hdslb.com/bfs/openplatform/202403/TXmgAvCC1710582339525.jpeg"
$ws.Range("C38").Value = "杭州·现世繁华-代号鸢only"
$ws.Range("D38").Value = "石祥路575号 杭州海外海纳川大酒店(万达广场渡驾桥地铁站店)"
$ws.Range("E38").Value = "2024.05.18 10:00-05.18 21:00"
$ws.Range("F38").Value = 585
$ws.Range("G38").Value = 76
$ws.Range("H38").Value = "https://show.bilibili.com/platform/detail.html?id=81905"
$ws.Range("I38").Value = "//i2.hdslb.com/bfs/openplatform/202402/m3upuV2F1708327958926.jpeg"
$cell = $ws.Range("B39")
$cell.NumberFormat = "@"
$cell.Value = "2024-05-25"
$cell.Style = "Normal"
$ws.Range("C39").Value = "杭州·原神X星铁X绝区零only"
$ws.Range("D39").Value = "望江东路333号 杭州瑞莱克斯大酒店"
$ws.Range("E39").Value = "2024.05.25 10:00-05.25 17:00"
$ws.Range("F39").Value = 105
$ws.Range("G39").Value = 60
$ws.Range("H39").Value = "https://show.bilibili.com/platform/detail.html?id=82754"
$ws.Range("I39").Value = "//i1.hdslb.com/bfs/openplatform/202403/qA0LNJuF1710234461030.jpeg"
$ws.Range("C40").Value = "杭州·第三届缘起cp展 我们二次元的情人节！"
$ws.Range("D40").Value = "黄姑山路51-4号 0101park"
$ws.Range("E40").Value = "2024.05.25 10:00-05.26 17:00"
$ws.Range("F40").Value = 46
$ws.Range("G40").Value = 65
$ws.Range("H40").Value = "https://show.bilibili.com/platform/detail.html?id=83336"
$ws.Range("I40").Value = "//i1.hdslb.com/bfs/openplatform/202403/D9t8ms7G1711350634757.png"
$cell = $ws.Range("B41")
$cell.NumberFormat = "@"
$cell.Value = "2024-06-01"
$cell.Style = "Normal"
$ws.Range("C41").Value = "杭州·造梦探险家——二次元同好会"
$ws.Range("D41").Value = "临平街道北沙西路156-1号 杭州临平遇上设计师酒店"
$ws.Range("E41").Value = "2024.06.01 10:00-06.01 16:00"
$ws.Range("F41").Value = 85
$ws.Range("G41").Value = 28
$ws.Range("H41").Value = "https://show.bilibili.com/platform/detail.html?id=82736"
$ws.Range("I41").Value = "//i1.hdslb.com/bfs/openplatform/202403/lqXD63661711623533572.png"
$cell = $ws.Range("B42")
$cell.NumberFormat = "@"
$cell.Value = "2024-06-08"
$cell.Style = "Normal"
$ws.Range("C42").Value = "杭州·第八届YH樱花动漫游戏文化节"
$ws.Range("D42").Value = "德胜东路2539号 梦马汽车小镇"
$ws.Range("E42").Value = "2024.06.08 10:00-06.10 17:00"
$ws.Range("F42").Value = 745
$ws.Range("G42").Value = 65
$ws.Range("H42").Value = "https://show.bilibili.com/platform/detail.html?id=82687"
$ws.Range("I42").Value = "//i2.hdslb.com/bfs/openplatform/202403/S5pnadXj1710210939138.png"
$cell = $ws.Range("B43")
$cell.NumberFormat = "@"
$cell.Value = "2024-06-09"
$cell.Style = "Normal"
$ws.Range("C43").Value = "杭州·第三届日夜国乙only"
$ws.Range("D43").Value = "创意路1号 中国智谷富春园区"
$ws.Range("E43").Value = "2024.06.09 10:00-06.09 23:00"
$ws.Range("F43").Value = 1425
$ws.Range("G43").Value = 58
$ws.Range("H43").Value = "https://show.bilibili.com/platform/detail.html?id=82618"
$ws.Range("I43").Value = "//i2.hdslb.com/bfs/openplatform/202403/fXRzYEFH1710124366279.png"
$cell = $ws.Range("B44")
$cell.NumberFormat = "@"
$cell.Value = "2024-06-14"
$cell.Style = "Normal"
$ws.Range("C44").Value = "杭州·苗阜王声 青曲社相声全国巡演"
$ws.Range("D44").Value = "湖墅南路138号 杭州浙话艺术剧院"
$ws.Range("E44").Value = "2024.06.14 19:30-06.14 22:00"
$ws.Range("F44").Value = 1
$ws.Range("G44").Value = 280
$ws.Range("H44").Value = "https://show.bilibili.com/platform/detail.html?id=83382"
$ws.Range("I44").Value = "//i1.hdslb.com/bfs/openplatform/202403/hUGL3xz01711346789039.jpeg"
$cell = $ws.Range("B45")
$cell.NumberFormat = "@"
$cell.Value = "2024-07-13"
$cell.Style = "Normal"
$ws.Range("C45").Value = "杭州·代号鸢only-广陵大学"
$ws.Range("D45").Value = "康候圣街99号 顺丰创新中心"
$ws.Range("E45").Value = "2024.07.13 09:00-07.13 18:00"
$ws.Range("F45").Value = 177
$ws.Range("G45").Value = 68
$ws.Range("H45").Value = "https://show.bilibili.com/platform/detail.html?id=83289"
$ws.Range("I45").Value = "//i0.hdslb.com/bfs/openplatform/202403/I3yffJ7Q1711344958258.png"
$cell = $ws.Range("B46")
$cell.NumberFormat = "@"
$cell.Value = "2024-07-20"
$cell.Style = "Normal"
$ws.Range("C46").Value = "杭州·次元幻想--二次元全女夜场"
$ws.Range("D46").Value = "保淑路2号 The Queen皇后"
$ws.Range("E46").Value = "2024.07.20 13:00-07.20 19:00"
$ws.Range("F46").Value = 468
$ws.Range("G46").Value = 158
$ws.Range("H46").Value = "https://show.bilibili.com/platform/detail.html?id=81808"
$ws.Range("I46").Value = "//i0.hdslb.com/bfs/openplatform/202402/sUUtSPh91707295826425.jpeg"
$cell = $ws.Range("B47")
$cell.NumberFormat = "@"
$cell.Value = "2024-08-03"
$cell.Style = "Normal"
$ws.Range("C47").Value = "杭州·梦漫星河动漫展"
$ws.Range("D47").Value = "阳城路雅澳杭州电商产业园西侧约200米 杭州大会展中心"
$ws.Range("E47").Value = "2024.08.03 10:00-08.04 17:00"
$ws.Range("F47").Value = 69
$ws.Range("G47").Value = 68
$ws.Range("H47").Value = "https://show.bilibili.com/platform/detail.html?id=82836"
$ws.Range("I47").Value = "//i0.hdslb.com/bfs/openplatform/202403/VFfQUJdD1711700169290.jpeg"
$ws.Range("F48").Value = 103
